$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate orders with updated distance/size codes.
# Distances: D51 -> D55, D64 -> D69, D80 -> D86
# Sizes:     S30 -> S31  (S20 and S25 stay unchanged)
#
# These substitutions touch every string that embeds the old tokens
# (Condition, Filename_Left, Filename_Right, Distance columns), so we
# perform a whole-sheet text replace for each token. xlWhole vs xlPart:
# use xlPart (2) since the tokens are embedded inside longer strings like
# "Face07_D51_S30" and "Face07_D51_S30_l.png".

$xlPart = 2
$xlByRows = 1

$ws.Cells.Replace("D51", "D55", $xlPart, $xlByRows, $false, $false, $true, $true)
$ws.Cells.Replace("D64", "D69", $xlPart, $xlByRows, $false, $false, $true, $true)
$ws.Cells.Replace("D80", "D86", $xlPart, $xlByRows, $false, $false, $true, $true)
$ws.Cells.Replace("S30", "S31", $xlPart, $xlByRows, $false, $false, $true, $true)
